# Update Gnai2-Tshr LR-pair data with recomputed TPM-based values.
# Source data (ligand/receptor average+total expression) changed for the
# "ECs", "MuSCs" and "Resolving-Mac" clusters (the "FAPs" cluster values are
# unchanged); every dependent column (specificity, edge weight, edge
# specificity) is updated to match.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    "G2" = 201.4397426666667
    "H2" = 604.3192280000001
    "I2" = 0.4833500233086392
    "J2" = 0.4833500233086393
    "M2" = 0.5550926666666666
    "N2" = 1.665278
    "O2" = 0.1208967663154349
    "P2" = 0.1208967663154349
    "Q2" = 111.8177239294871
    "R2" = 1006.359515365384
    "S2" = 0.05843545481650456
    "T2" = 0.05843545481650456
    "G3" = 201.4397426666667
    "H3" = 604.3192280000001
    "I3" = 0.4833500233086392
    "J3" = 0.4833500233086393
    "O3" = 0.7377399926530269
    "P3" = 0.7377399926530268
    "Q3" = 682.3375789471899
    "R3" = 6141.038210524709
    "S3" = 0.3565866426445559
    "T3" = 0.3565866426445559
    "G4" = 201.4397426666667
    "H4" = 604.3192280000001
    "I4" = 0.4833500233086392
    "J4" = 0.4833500233086393
    "M4" = 0.5311786666666667
    "N4" = 1.593536
    "O4" = 0.1156884012202364
    "P4" = 0.1156884012202364
    "Q4" = 107.0004939233565
    "R4" = 963.0044453102081
    "S4" = 0.05591799142634048
    "T4" = 0.05591799142634048
    "G5" = 201.4397426666667
    "H5" = 604.3192280000001
    "I5" = 0.4833500233086392
    "J5" = 0.4833500233086393
    "M5" = 0.117885
    "N5" = 0.353655
    "O5" = 0.02567483981130185
    "P5" = 0.02567483981130185
    "Q5" = 23.74672406426
    "R5" = 213.72051657834
    "S5" = 0.01240993442123833
    "T5" = 0.01240993442123833
    "I6" = 0.1569674599353791
    "J6" = 0.1569674599353792
    "M6" = 0.5550926666666666
    "N6" = 1.665278
    "O6" = 0.1208967663154349
    "P6" = 0.1208967663154349
    "Q6" = 36.31269939912578
    "R6" = 326.8142945921319
    "S6" = 0.01897685832293492
    "T6" = 0.01897685832293492
    "I7" = 0.1569674599353791
    "J7" = 0.1569674599353792
    "O7" = 0.7377399926530269
    "P7" = 0.7377399926530268
    "S7" = 0.1158011727394909
    "T7" = 0.1158011727394909
    "I8" = 0.1569674599353791
    "J8" = 0.1569674599353792
    "M8" = 0.5311786666666667
    "N8" = 1.593536
    "O8" = 0.1156884012202364
    "P8" = 0.1156884012202364
    "Q8" = 34.74830854048712
    "R8" = 312.734776864384
    "S8" = 0.01815931448352553
    "T8" = 0.01815931448352553
    "I9" = 0.1569674599353791
    "J9" = 0.1569674599353792
    "M9" = 0.117885
    "N9" = 0.353655
    "O9" = 0.02567483981130185
    "P9" = 0.02567483981130185
    "Q9" = 7.711726033730001
    "R9" = 69.40553430356999
    "S9" = 0.004030114389427801
    "T9" = 0.004030114389427801
    "G10" = 60.43484133333334
    "H10" = 181.304524
    "I10" = 0.1450120099461104
    "J10" = 0.1450120099461104
    "M10" = 0.5550926666666666
    "N10" = 1.665278
    "O10" = 0.1208967663154349
    "P10" = 0.1208967663154349
    "Q10" = 33.54693723529689
    "R10" = 301.922435117672
    "S10" = 0.01753148307938642
    "T10" = 0.01753148307938642
    "G11" = 60.43484133333334
    "H11" = 181.304524
    "I11" = 0.1450120099461104
    "J11" = 0.1450120099461104
    "O11" = 0.7377399926530269
    "P11" = 0.7377399926530268
    "Q11" = 204.7111596428182
    "R11" = 1842.400436785364
    "S11" = 0.1069811591522441
    "T11" = 0.1069811591522441
    "G12" = 60.43484133333334
    "H12" = 181.304524
    "I12" = 0.1450120099461104
    "J12" = 0.1450120099461104
    "M12" = 0.5311786666666667
    "N12" = 1.593536
    "O12" = 0.1156884012202364
    "P12" = 0.1156884012202364
    "Q12" = 32.10169843965156
    "R12" = 288.915285956864
    "S12" = 0.01677620758839853
    "T12" = 0.01677620758839853
    "G13" = 60.43484133333334
    "H13" = 181.304524
    "I13" = 0.1450120099461104
    "J13" = 0.1450120099461104
    "M13" = 0.117885
    "N13" = 0.353655
    "O13" = 0.02567483981130185
    "P13" = 0.02567483981130185
    "Q13" = 7.124361270580001
    "R13" = 64.11925143522001
    "S13" = 0.003723160126081294
    "T13" = 0.003723160126081295
    "G14" = 89.46554166666668
    "H14" = 268.396625
    "I14" = 0.2146705068098712
    "J14" = 0.2146705068098712
    "M14" = 0.5550926666666666
    "N14" = 1.665278
    "O14" = 0.1208967663154349
    "P14" = 0.1208967663154349
    "Q14" = 49.66166609852777
    "R14" = 446.95499488675
    "S14" = 0.02595297009660897
    "T14" = 0.02595297009660897
    "G15" = 89.46554166666668
    "H15" = 268.396625
    "I15" = 0.2146705068098712
    "J15" = 0.2146705068098712
    "O15" = 0.7377399926530269
    "P15" = 0.7377399926530268
    "Q15" = 303.0469573278195
    "R15" = 2727.422615950375
    "S15" = 0.158371018116736
    "T15" = 0.158371018116736
    "G16" = 89.46554166666668
    "H16" = 268.396625
    "I16" = 0.2146705068098712
    "J16" = 0.2146705068098712
    "M16" = 0.5311786666666667
    "N16" = 1.593536
    "O16" = 0.1156884012202364
    "P16" = 0.1156884012202364
    "Q16" = 47.52218713511112
    "R16" = 427.699684216
    "S16" = 0.02483488772197187
    "T16" = 0.02483488772197187
    "G17" = 89.46554166666668
    "H17" = 268.396625
    "I17" = 0.2146705068098712
    "J17" = 0.2146705068098712
    "M17" = 0.117885
    "N17" = 0.353655
    "O17" = 0.02567483981130185
    "P17" = 0.02567483981130185
    "Q17" = 10.546645379375
    "R17" = 94.91980841437501
    "S17" = 0.005511630874554427
    "T17" = 0.005511630874554427
}

foreach ($cellRef in $newValues.Keys) {
    $ws.Range($cellRef).Value2 = $newValues[$cellRef]
}
